$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Replace "Environmental Perturbation" with "Environmental_Perturbation" (column D, rows 2-27) ---
$ws.Range("D2:D27").Value = "Environmental_Perturbation"

# --- Replace "KN99 alpha" with "KN99_alpha" (column F, rows 2-28,30,32,34,36) ---
$ws.Range("F2:F28").Value = "KN99_alpha"
$ws.Range("F30").Value = "KN99_alpha"
$ws.Range("F32").Value = "KN99_alpha"
$ws.Range("F34").Value = "KN99_alpha"
$ws.Range("F36").Value = "KN99_alpha"

# --- Replace "Time Course" with "Timecourse" (column D, rows 28-37) ---
$ws.Range("D28:D37").Value = "Timecourse"

# --- Update sheet view: topLeftCell and selection ---
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("J21").Select()
